# Daily attendance processing - 2025-10-15 23:20:49
# Reorders the 'Recorded By' (column G) names in specific rows so that the
# last-listed recorder becomes the first-listed recorder (cosmetic re-ordering
# of the comma-separated attendee/recorder names already present in the cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = 'backup@backdoor.com, System, system'
    3 = 'dnasr281@gmail.com, System'
    4 = 'System, backup@backdoor.com'
    5 = 'backup@backdoor.com, System'
    6 = 'dnasr281@gmail.com, System'
    7 = 'admin@admin.com, System'
    10 = 'dnasr281@gmail.com, System'
    12 = 'dnasr281@gmail.com, System'
    13 = 'dnasr281@gmail.com, System'
    14 = 'dnasr281@gmail.com, System'
    15 = 'dnasr281@gmail.com, System'
    19 = 'dnasr281@gmail.com, System'
    29 = 'backup@backdoor.com, System, system'
    30 = 'dnasr281@gmail.com, System'
    31 = 'System, backup@backdoor.com'
    32 = 'backup@backdoor.com, System'
    33 = 'dnasr281@gmail.com, System'
    34 = 'admin@admin.com, System'
    37 = 'dnasr281@gmail.com, System'
    39 = 'dnasr281@gmail.com, System'
    40 = 'dnasr281@gmail.com, System'
    41 = 'dnasr281@gmail.com, System'
    42 = 'dnasr281@gmail.com, System'
    46 = 'dnasr281@gmail.com, System'
    56 = 'backup@backdoor.com, System, system'
    57 = 'dnasr281@gmail.com, System'
    58 = 'System, backup@backdoor.com'
    59 = 'backup@backdoor.com, System'
    60 = 'dnasr281@gmail.com, System'
    61 = 'admin@admin.com, System'
    64 = 'dnasr281@gmail.com, System'
    66 = 'dnasr281@gmail.com, System'
    67 = 'dnasr281@gmail.com, System'
    68 = 'dnasr281@gmail.com, System'
    69 = 'dnasr281@gmail.com, System'
    73 = 'dnasr281@gmail.com, System'
    84 = 'backup@backdoor.com, System'
    85 = 'backup@backdoor.com, System'
    86 = 'dnasr281@gmail.com, System'
    87 = 'dnasr281@gmail.com, System'
    88 = 'dnasr281@gmail.com, System'
    89 = 'dnasr281@gmail.com, System'
    93 = 'dnasr281@gmail.com, System'
    95 = 'dnasr281@gmail.com, System'
    110 = 'backup@backdoor.com, System'
    111 = 'backup@backdoor.com, System'
    112 = 'dnasr281@gmail.com, System'
    113 = 'dnasr281@gmail.com, System'
    114 = 'dnasr281@gmail.com, System'
    115 = 'dnasr281@gmail.com, System'
    119 = 'dnasr281@gmail.com, System'
    121 = 'dnasr281@gmail.com, System'
    136 = 'backup@backdoor.com, System'
    137 = 'backup@backdoor.com, System'
    138 = 'dnasr281@gmail.com, System'
    139 = 'dnasr281@gmail.com, System'
    140 = 'dnasr281@gmail.com, System'
    141 = 'dnasr281@gmail.com, System'
    145 = 'dnasr281@gmail.com, System'
    147 = 'dnasr281@gmail.com, System'
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}

